# Auto-generated: apply scheduled-runner market-data refresh to Sheets/Behemoth_Profits.xlsx
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) per refreshed market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 919.64703
$ws.Range("I28").Value = 718.46155
$ws.Range("K28").Value = 718.46155
$ws.Range("M28").Value = -233.46155
$ws.Range("H69").Value = 28335.666
$ws.Range("J69").Value = 38507.5
$ws.Range("L69").Value = 115522.5
$ws.Range("N69").Value = -117270.5
$ws.Range("H72").Value = 28335.666
$ws.Range("J72").Value = 38507.5
$ws.Range("L72").Value = 346567.5
$ws.Range("N72").Value = -355303.5
$ws.Range("H111").Value = 3019.6
$ws.Range("I111").Value = 2910.7778
$ws.Range("K111").Value = 8732.3334
$ws.Range("M111").Value = -5665.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 4062666.2
$ws.Range("I16").Value = 6091249.5
$ws.Range("J16").Value = 5499.5
$ws.Range("K16").Value = 6091249.5
$ws.Range("L16").Value = 5499.5
$ws.Range("M16").Value = -6090962.5
$ws.Range("N16").Value = -6073.5
$ws.Range("H45").Value = 2505
$ws.Range("J45").Value = 2902.2
$ws.Range("L45").Value = 2902.2
$ws.Range("N45").Value = -3656.2
$ws.Range("H46").Value = 4901
$ws.Range("J46").Value = 4901
$ws.Range("L46").Value = 4901
$ws.Range("N46").Value = -5539
$ws.Range("H74").Value = 10010086
$ws.Range("I74").Value = 13891899
$ws.Range("J74").Value = 28281.715
$ws.Range("K74").Value = 13891899
$ws.Range("L74").Value = 28281.715
$ws.Range("M74").Value = -13891025
$ws.Range("N74").Value = -30029.715
$ws.Range("H77").Value = 10010086
$ws.Range("I77").Value = 13891899
$ws.Range("J77").Value = 28281.715
$ws.Range("K77").Value = 69459495
$ws.Range("L77").Value = 141408.575
$ws.Range("M77").Value = -69455127
$ws.Range("N77").Value = -150144.575
$ws.Range("H132").Value = 11019.147
$ws.Range("I132").Value = 8780.956
$ws.Range("J132").Value = 15699
$ws.Range("K132").Value = 26342.868
$ws.Range("L132").Value = 47097
$ws.Range("M132").Value = -23812.868
$ws.Range("N132").Value = -52157

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 389.9091
$ws.Range("I22").Value = 389.9091
$ws.Range("K22").Value = 389.9091
$ws.Range("M22").Value = -216.9091
$ws.Range("H86").Value = 2786.6667
$ws.Range("I86").Value = 2672.5
$ws.Range("J86").Value = 3700
$ws.Range("K86").Value = 2672.5
$ws.Range("L86").Value = 3700
$ws.Range("M86").Value = -1549.5
$ws.Range("N86").Value = -5946
$ws.Range("H89").Value = 2786.6667
$ws.Range("I89").Value = 2672.5
$ws.Range("J89").Value = 3700
$ws.Range("K89").Value = 13362.5
$ws.Range("L89").Value = 18500
$ws.Range("M89").Value = -7746.5
$ws.Range("N89").Value = -29732
$ws.Range("H105").Value = 2732.6667
$ws.Range("I105").Value = 1599.5
$ws.Range("K105").Value = 1599.5
$ws.Range("M105").Value = 147.5
$ws.Range("H107").Value = 1405.2858
$ws.Range("I107").Value = 1385
$ws.Range("J107").Value = 1669
$ws.Range("K107").Value = 1385
$ws.Range("L107").Value = 1669
$ws.Range("M107").Value = 535
$ws.Range("N107").Value = -5509
$ws.Range("H114").Value = 107000
$ws.Range("J114").Value = 107000
$ws.Range("L114").Value = 107000
$ws.Range("N114").Value = -115678
$ws.Range("H133").Value = 58000
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("H134").Value = 53715.57
$ws.Range("I134").Value = 2865.6667
$ws.Range("K134").Value = 8597.000100000001
$ws.Range("M134").Value = -6062.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2937
$ws.Range("I107").Value = 2937
$ws.Range("K107").Value = 2937
$ws.Range("M107").Value = -1017
$ws.Range("H125").Value = 252650.75
$ws.Range("J125").Value = 252650.75
$ws.Range("L125").Value = 252650.75
$ws.Range("N125").Value = -257570.75
$ws.Range("H132").Value = 3632.4
$ws.Range("I132").Value = 3165.5
$ws.Range("J132").Value = 5500
$ws.Range("K132").Value = 9496.5
$ws.Range("L132").Value = 16500
$ws.Range("M132").Value = -6966.5
$ws.Range("N132").Value = -21560

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 9998.75
$ws.Range("I56").Value = 9998.75
$ws.Range("K56").Value = 9998.75
$ws.Range("M56").Value = -9468.75
$ws.Range("H140").Value = 179317.94
$ws.Range("I140").Value = 190244.06
$ws.Range("K140").Value = 570732.1799999999
$ws.Range("M140").Value = -565552.1799999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 10000
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 10000
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 10000
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -10340
$ws.Range("H97").Value = 897.4167
$ws.Range("I97").Value = 997.2632
$ws.Range("K97").Value = 997.2632
$ws.Range("M97").Value = -501.2632
$ws.Range("H113").Value = 4060.7144
$ws.Range("J113").Value = 4460.25
$ws.Range("L113").Value = 4460.25
$ws.Range("N113").Value = -8800.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 10003
$ws.Range("I22").Value = 10003
$ws.Range("K22").Value = 10003
$ws.Range("M22").Value = -9708
$ws.Range("H27").Value = 10003
$ws.Range("I27").Value = 10003
$ws.Range("K27").Value = 10003
$ws.Range("M27").Value = -9896
$ws.Range("H32").Value = 8933.333000000001
$ws.Range("I32").Value = 13000
$ws.Range("J32").Value = 800
$ws.Range("K32").Value = 13000
$ws.Range("L32").Value = 800
$ws.Range("M32").Value = -12683
$ws.Range("N32").Value = -1434
$ws.Range("H40").Value = 3442.8667
$ws.Range("I40").Value = 2812.4783
$ws.Range("K40").Value = 2812.4783
$ws.Range("M40").Value = -2676.4783
$ws.Range("H122").Value = 4966.7144
$ws.Range("I122").Value = 3737.5625
$ws.Range("K122").Value = 11212.6875
$ws.Range("M122").Value = -8762.6875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4262177
$ws.Range("J62").Value = 22225786
$ws.Range("L62").Value = 22225786
$ws.Range("N62").Value = -22227034
$ws.Range("H65").Value = 4262177
$ws.Range("J65").Value = 22225786
$ws.Range("L65").Value = 111128930
$ws.Range("N65").Value = -111135170
$ws.Range("H126").Value = 3467.75
$ws.Range("I126").Value = 2686.25
$ws.Range("K126").Value = 8058.75
$ws.Range("M126").Value = -5588.75
$ws.Range("H132").Value = 17157.615
$ws.Range("I132").Value = 2391
$ws.Range("J132").Value = 26386.75
$ws.Range("K132").Value = 7173
$ws.Range("L132").Value = 79160.25
$ws.Range("M132").Value = -4643
$ws.Range("N132").Value = -84220.25

